$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.331.12'
$ws.Range("E2").Value = '  +3.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.089.83'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.33'
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.75'
$ws.Range("E6").Value = '  -2.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.373'
$ws.Range("E7").Value = '  -4.45%  '
$ws.Range("E8").Value = '  +8.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.086.57'
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.667'
$ws.Range("E11").Value = '  +17.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.189'
$ws.Range("E12").Value = '  +5.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '91.142.61'
$ws.Range("E14").Value = '  +3.70%  '
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.91'
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.658.82'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.078.21'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  +3.77%  '
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.78'
$ws.Range("E21").Value = '  +3.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '434.69'
$ws.Range("E22").Value = '  +2.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.46'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.12'
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.62'
$ws.Range("E25").Value = '  +2.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.84'
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.79'
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.257.50'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.165'
$ws.Range("E30").Value = '  +6.49%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.60'
$ws.Range("E32").Value = '  +5.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.83'
$ws.Range("E33").Value = '  -6.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '514.82'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.94'
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.140'
$ws.Range("E36").Value = '  -6.80%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.04'
$ws.Range("E39").Value = '  +3.32%  '
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +3.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.368'
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("E45").Value = '  +1.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0716'
$ws.Range("E46").Value = '  +10.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.75'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.92'
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.21'
$ws.Range("E49").Value = '  +6.65%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000260'
$ws.Range("E50").Value = '  +9.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '165.24'
$ws.Range("E51").Value = '  +1.57%  '
